# The commit swaps the colour scheme carried by the deck's theme part
# (ppt/theme/theme1.xml, the theme used by the slide master / all slides)
# from the "Integral" / "Red Violet" palette over to the stock Office
# "Office Theme" palette (the palette that used to only live in
# ppt/theme/theme2.xml, which the Notes Master referenced). The theme's
# font scheme and format scheme (fills/lines/effects) are identical
# between the two themes already, so only the 12 colour-scheme slots
# need to change.
#
# PowerPoint's object model exposes those 12 slots through
# Theme.ThemeColorScheme.Colors(index).RGB (1-based, in the fixed order
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) - exactly the values
# stored in <a:clrScheme> for the theme.

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> (scheme slot, target "Office Theme" RGB)
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
